# Update the BIA variable-name column (B) to be prefixed with "BIA_" for
# all rows below the Sex/Age rows (rows 5-23), to reflect the updated
# data dictionary ("all_data_dicts_Jan_2018").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$varNames = @{
    5  = "BIA_Frame"
    6  = "BIA_Activity_Level"
    7  = "BIA_Height"
    8  = "BIA_Weight"
    10 = "BIA_Fat"
    11 = "BIA_LDM"
    12 = "BIA_ICW"
    13 = "BIA_BMI"
    14 = "BIA_ECW"
    15 = "BIA_FFM"
    16 = "BIA_TBW"
    17 = "BIA_SMM"
    18 = "BIA_FMI"
    19 = "BIA_FFMI"
    20 = "BIA_BMC"
    21 = "BIA_LST"
    22 = "BIA_BMR"
    23 = "BIA_DEE"
}

foreach ($row in $varNames.Keys) {
    $ws.Cells.Item($row, 2).Value = $varNames[$row]
}

# Move the active selection to B10, matching the saved selection state.
$ws.Range("B10").Select()
